$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add "Y" values in columns A and B (columns C..AD already contain "Y")
$ws.Range("A4").Value = "Y"
$ws.Range("B4").Value = "Y"

# Row 7: new row of "Y" values across columns A..J and L..Q (K left blank)
$ws.Range("A7").Value = "Y"
$ws.Range("B7").Value = "Y"
$ws.Range("C7").Value = "Y"
$ws.Range("D7").Value = "Y"
$ws.Range("E7").Value = "Y"
$ws.Range("F7").Value = "Y"
$ws.Range("G7").Value = "Y"
$ws.Range("H7").Value = "Y"
$ws.Range("I7").Value = "Y"
$ws.Range("J7").Value = "Y"
$ws.Range("L7").Value = "Y"
$ws.Range("M7").Value = "Y"
$ws.Range("N7").Value = "Y"
$ws.Range("O7").Value = "Y"
$ws.Range("P7").Value = "Y"
$ws.Range("Q7").Value = "Y"

# Update the active cell selection to H5
$ws.Range("H5").Select()
